# Insert two new weekly price rows for "Feria Lagunitas de Puerto Montt - Lechuga"
# right before the existing row that used to be row 697 (Mercado ID column = A,
# date serial 44467 / 2021-09-28), pushing all subsequent rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at 697-698 (everything from the old row 697 onward moves
# down two rows, so old row 775 becomes new row 777).
$ws.Range("A697:A698").EntireRow.Insert()

# --- New row 697: Lechuga / Escarola / Primera, Región de Coquimbo ---
$ws.Cells.Item(697, 1).Value2  = 4
$ws.Cells.Item(697, 2).Value   = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(697, 3).Value   = "Los Lagos"
$ws.Cells.Item(697, 4).Value2  = 44946
$ws.Cells.Item(697, 5).Value2  = 10
$ws.Cells.Item(697, 6).Value2  = 100112033
$ws.Cells.Item(697, 7).Value   = "Lechuga"
$ws.Cells.Item(697, 8).Value   = "Escarola"
$ws.Cells.Item(697, 9).Value   = "Primera"
$ws.Cells.Item(697, 10).Value2 = 300
$ws.Cells.Item(697, 11).Value2 = 15000
$ws.Cells.Item(697, 12).Value2 = 15000
$ws.Cells.Item(697, 13).Value2 = 15000
$ws.Cells.Item(697, 14).Value  = "`$/caja 15 unidades"
$ws.Cells.Item(697, 15).Value  = "Región de Coquimbo"
$ws.Cells.Item(697, 16).Value2 = 1000
$ws.Cells.Item(697, 17).Value2 = 15
$ws.Cells.Item(697, 18).Value  = "Hortaliza"

# --- New row 698: Lechuga / Escarola / Segunda, Región de Coquimbo ---
$ws.Cells.Item(698, 1).Value2  = 4
$ws.Cells.Item(698, 2).Value   = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(698, 3).Value   = "Los Lagos"
$ws.Cells.Item(698, 4).Value2  = 44946
$ws.Cells.Item(698, 5).Value2  = 10
$ws.Cells.Item(698, 6).Value2  = 100112033
$ws.Cells.Item(698, 7).Value   = "Lechuga"
$ws.Cells.Item(698, 8).Value   = "Escarola"
$ws.Cells.Item(698, 9).Value   = "Segunda"
$ws.Cells.Item(698, 10).Value2 = 300
$ws.Cells.Item(698, 11).Value2 = 13000
$ws.Cells.Item(698, 12).Value2 = 13000
$ws.Cells.Item(698, 13).Value2 = 13000
$ws.Cells.Item(698, 14).Value  = "`$/caja 18 unidades"
$ws.Cells.Item(698, 15).Value  = "Región de Coquimbo"
$ws.Cells.Item(698, 16).Value2 = 722
$ws.Cells.Item(698, 17).Value2 = 18
$ws.Cells.Item(698, 18).Value  = "Hortaliza"
